$d = $word.ActiveDocument

# Remove the whole "USCOTS 2013 Notes" paragraph (paragraph 10), including its
# paragraph mark, collapsing the surrounding blank paragraphs.
$titlePara = $d.Paragraphs.Item(10)
$titlePara.Range.Delete()

# The "_GoBack" bookmark previously sat in the paragraph right before this
# block (now paragraph 7); re-adding it under the same name moves it to the
# new home paragraph (now paragraph 10, the last of the plain blank
# paragraphs that used to be paragraph 11).
$newBookmarkHome = $d.Paragraphs.Item(10)
$d.Bookmarks.Add("_GoBack", $newBookmarkHome.Range)
